$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Memoria por correspondencia"
$ws.Range("B18").Value = "Emma Reyes"

# C18 stays empty (matching the pattern of the other rows), but we still
# need a cell record to exist for it, so nudge a formatting property that
# is already at its default value - this materializes the cell without
# changing its appearance or introducing any new style.
$ws.Range("C18").Borders.LineStyle = -4142
